# Generate Report for Archive
#
# The localization run picked up new status for bc48b532-...md (now
# "In Translation", same as 356a4b89-...md) and re-ordered the two
# "Ready for handoff" / "In Translation" rows (9975d285-...md and
# bc48b532-...md) in each of the three sheets: Overview, zh-cn, de-de.
#
# Logical effect (verified against the target OOXML diff):
#   - 356a4b89-...md : Status "Ready for handoff" -> "In Translation"
#   - bc48b532-...md : Status "Ready for handoff" -> "In Translation";
#                       row moves up to the old 9975d285 row slot
#   - 9975d285-...md : Status stays "Ready for handoff"; row moves down
#                       to the old bc48b532 row slot
#
# Concretely this swaps the data of row 4 / row 5 (keeping each file's
# own dates / xlf hash filenames attached to it) and updates the status
# text, in all three worksheets.

$wb = $excel.ActiveWorkbook

$bc48Url9975 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1952973248f04b79b17c9dbe07464e2ea7b135da/e2e/9975d285-3142-4e0b-82a7-6274ead42379.md"
$bc48UrlBc48 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae3d4e80bdfb593cba55ecb8d8c80a382e123779/e2e/bc48b532-de57-482e-8253-c71468617fb9.md"
$url677a     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/67009a62cfb0e1019bb316c3868655d616089b43/e2e/677a229b-e090-4c61-a718-ff5724991837.md"
$url356a     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ae3d4e80bdfb593cba55ecb8d8c80a382e123779/e2e/356a4b89-232d-4b0c-9593-dbe600d1b493.md"

# ---------------------------------------------------------------------
# Sheet "Overview" : columns A=File Name, B=Path And Name (hyperlink),
#                     C=Extension, D=Publish URL, E=zh-cn, F=de-de,
#                     G=Latest HO Xliff Generate Date
# Row4 becomes bc48b532's data, Row5 becomes 9975d285's data.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "bc48b532-de57-482e-8253-c71468617fb9.md"
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"
$wsOverview.Range("G4").Value = "2016-08-23 22:38:34"

$wsOverview.Range("A5").Value = "9975d285-3142-4e0b-82a7-6274ead42379.md"
$wsOverview.Range("E5").Value = "Ready for handoff"
$wsOverview.Range("F5").Value = "Ready for handoff"
$wsOverview.Range("G5").Value = "2016-08-23 22:37:14"

# Status for 356a4b89 (row 3) also flips to "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Rebuild the B-column hyperlinks (display text + target) so B4 now
# points at bc48b532 and B5 now points at 9975d285.
$wsOverview.UsedRange.Hyperlinks.Delete()
$ovLinks = $wsOverview.Hyperlinks
$ovLinks.Add($wsOverview.Range("B2"), $url677a, "", "", "e2e\677a229b-e090-4c61-a718-ff5724991837.md")
$ovLinks.Add($wsOverview.Range("B3"), $url356a, "", "", "e2e\356a4b89-232d-4b0c-9593-dbe600d1b493.md")
$ovLinks.Add($wsOverview.Range("B4"), $bc48UrlBc48, "", "", "e2e\bc48b532-de57-482e-8253-c71468617fb9.md")
$ovLinks.Add($wsOverview.Range("B5"), $bc48Url9975, "", "", "e2e\9975d285-3142-4e0b-82a7-6274ead42379.md")

# ---------------------------------------------------------------------
# Sheet "zh-cn" : columns A=Source File Name, C=Status, G=Latest
#                 Handoff File, H=Latest Handoff Datetime (+ hyperlink
#                 on column A)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# 356a4b89 (row 3): status flips, own hash/date are unchanged.
$wsZh.Range("C3").Value = "In Translation"

# Row 4 becomes bc48b532's own data (status "In Translation").
$wsZh.Range("A4").Value = "bc48b532-de57-482e-8253-c71468617fb9.md"
$wsZh.Range("C4").Value = "In Translation"
$wsZh.Range("G4").Value = "bc48b532-de57-482e-8253-c71468617fb9.a1c827bff2686d01c93ccc970d489682cf4503fa.zh-cn.xlf"
$wsZh.Range("H4").Value = "2016-08-23 22:38:29"

# Row 5 becomes 9975d285's own data (status "Ready for handoff").
$wsZh.Range("A5").Value = "9975d285-3142-4e0b-82a7-6274ead42379.md"
$wsZh.Range("C5").Value = "Ready for handoff"
$wsZh.Range("G5").Value = "9975d285-3142-4e0b-82a7-6274ead42379.3ca77b64f2b02cafb0fa5456deac5076d84d6329.zh-cn.xlf"
$wsZh.Range("H5").Value = "2016-08-23 22:37:09"

$wsZh.UsedRange.Hyperlinks.Delete()
$zhLinks = $wsZh.Hyperlinks
$zhLinks.Add($wsZh.Range("A2"), $url677a, "", "", "677a229b-e090-4c61-a718-ff5724991837.md")
$zhLinks.Add($wsZh.Range("I2"), $url677a, "", "", "677a229b-e090-4c61-a718-ff5724991837.md")
$zhLinks.Add($wsZh.Range("A3"), $url356a, "", "", "356a4b89-232d-4b0c-9593-dbe600d1b493.md")
$zhLinks.Add($wsZh.Range("A4"), $bc48UrlBc48, "", "", "bc48b532-de57-482e-8253-c71468617fb9.md")
$zhLinks.Add($wsZh.Range("A5"), $bc48Url9975, "", "", "9975d285-3142-4e0b-82a7-6274ead42379.md")

# ---------------------------------------------------------------------
# Sheet "de-de" : same layout/columns as zh-cn.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# 356a4b89 (row 3): status flips, own hash/date are unchanged.
$wsDe.Range("C3").Value = "In Translation"

# Row 4 becomes bc48b532's own data (status "In Translation").
$wsDe.Range("A4").Value = "bc48b532-de57-482e-8253-c71468617fb9.md"
$wsDe.Range("C4").Value = "In Translation"
$wsDe.Range("G4").Value = "bc48b532-de57-482e-8253-c71468617fb9.a1c827bff2686d01c93ccc970d489682cf4503fa.de-de.xlf"
$wsDe.Range("H4").Value = "2016-08-23 22:38:34"

# Row 5 becomes 9975d285's own data (status "Ready for handoff").
$wsDe.Range("A5").Value = "9975d285-3142-4e0b-82a7-6274ead42379.md"
$wsDe.Range("C5").Value = "Ready for handoff"
$wsDe.Range("G5").Value = "9975d285-3142-4e0b-82a7-6274ead42379.3ca77b64f2b02cafb0fa5456deac5076d84d6329.de-de.xlf"
$wsDe.Range("H5").Value = "2016-08-23 22:37:14"

$wsDe.UsedRange.Hyperlinks.Delete()
$deLinks = $wsDe.Hyperlinks
$deLinks.Add($wsDe.Range("A2"), $url677a, "", "", "677a229b-e090-4c61-a718-ff5724991837.md")
$deLinks.Add($wsDe.Range("I2"), $url677a, "", "", "677a229b-e090-4c61-a718-ff5724991837.md")
$deLinks.Add($wsDe.Range("A3"), $url356a, "", "", "356a4b89-232d-4b0c-9593-dbe600d1b493.md")
$deLinks.Add($wsDe.Range("A4"), $bc48UrlBc48, "", "", "bc48b532-de57-482e-8253-c71468617fb9.md")
$deLinks.Add($wsDe.Range("A5"), $bc48Url9975, "", "", "9975d285-3142-4e0b-82a7-6274ead42379.md")

$wsOverview.Select()
